$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append " (23 base core)" to the three spacer header labels in row 1
$ws.Range("B1").Value = 'Spacer "A" (23 base core)'
$ws.Range("C1").Value = 'Spacer "B" (23 base core)'
$ws.Range("D1").Value = 'Spacer "C" (23 base core)'

# Update the view: scroll so column C is the left-most visible column,
# and select E5 (mirrors the saved sheetView state in the target workbook)
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("E5").Select() | Out-Null
